$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 57.444443
$ws.Range("I4").Value = 68.14286
$ws.Range("K4").Value = 68.14286
$ws.Range("M4").Value = 45.85714
$ws.Range("H15").Value = 1671.119
$ws.Range("I15").Value = 1671.119
$ws.Range("K15").Value = 5013.357
$ws.Range("M15").Value = -4844.357
$ws.Range("H43").Value = 4366.1333
$ws.Range("I43").Value = 3548.3333
$ws.Range("J43").Value = 4911.3335
$ws.Range("K43").Value = 3548.3333
$ws.Range("L43").Value = 4911.3335
$ws.Range("M43").Value = -3479.3333
$ws.Range("N43").Value = -5049.3335
$ws.Range("H132").Value = 1741.7115
$ws.Range("I132").Value = 1416.766
$ws.Range("K132").Value = 4250.298000000001
$ws.Range("M132").Value = -1720.298000000001
$ws.Range("H137").Value = 4275955.5
$ws.Range("I137").Value = 2686.9
$ws.Range("K137").Value = 8060.700000000001
$ws.Range("M137").Value = -5510.700000000001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17490010
$ws.Range("I32").Value = 17721864
$ws.Range("K32").Value = 17721864
$ws.Range("M32").Value = -17721577
$ws.Range("H45").Value = 5219.3335
$ws.Range("I45").Value = 4663.2
$ws.Range("K45").Value = 4663.2
$ws.Range("M45").Value = -4286.2
$ws.Range("H61").Value = 3662
$ws.Range("I61").Value = 3270.7144
$ws.Range("K61").Value = 3270.7144
$ws.Range("M61").Value = -3058.7144
$ws.Range("H74").Value = 2762.5
$ws.Range("I74").Value = 2699.6316
$ws.Range("K74").Value = 2699.6316
$ws.Range("M74").Value = -1825.6316
$ws.Range("H77").Value = 2762.5
$ws.Range("I77").Value = 2699.6316
$ws.Range("K77").Value = 13498.158
$ws.Range("M77").Value = -9130.158000000001
$ws.Range("H122").Value = 3238.125
$ws.Range("I122").Value = 3150.8333
$ws.Range("K122").Value = 9452.499899999999
$ws.Range("M122").Value = -7002.499899999999
$ws.Range("H136").Value = 3662
$ws.Range("I136").Value = 3270.7144
$ws.Range("K136").Value = 9812.143199999999
$ws.Range("M136").Value = -7262.143199999999

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50181
$ws.Range("J62").Value = 50181
$ws.Range("L62").Value = 50181
$ws.Range("N62").Value = -51553
$ws.Range("H65").Value = 50181
$ws.Range("J65").Value = 50181
$ws.Range("L65").Value = 150543
$ws.Range("N65").Value = -157407
$ws.Range("H94").Value = 1419.2354
$ws.Range("I94").Value = 885.25
$ws.Range("J94").Value = 1893.8889
$ws.Range("K94").Value = 885.25
$ws.Range("L94").Value = 1893.8889
$ws.Range("M94").Value = -434.25
$ws.Range("N94").Value = -2795.8889
$ws.Range("H105").Value = 2030.6316
$ws.Range("I105").Value = 1253.3636
$ws.Range("K105").Value = 1253.3636
$ws.Range("M105").Value = 493.6364000000001
$ws.Range("H134").Value = 1589116.2
$ws.Range("I134").Value = 1787236.2
$ws.Range("K134").Value = 5361708.6
$ws.Range("M134").Value = -5359173.6

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5666.967
$ws.Range("I31").Value = 2099.25
$ws.Range("J31").Value = 6215.846
$ws.Range("K31").Value = 2099.25
$ws.Range("L31").Value = 6215.846
$ws.Range("M31").Value = -1804.25
$ws.Range("N31").Value = -6805.846
$ws.Range("H34").Value = 5666.967
$ws.Range("I34").Value = 2099.25
$ws.Range("J34").Value = 6215.846
$ws.Range("K34").Value = 2099.25
$ws.Range("L34").Value = 6215.846
$ws.Range("M34").Value = -1897.25
$ws.Range("N34").Value = -6619.846
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3754.7693
$ws.Range("I132").Value = 3892
$ws.Range("K132").Value = 11676
$ws.Range("M132").Value = -9146
$ws.Range("H134").Value = 1600.8125
$ws.Range("I134").Value = 1439.9333
$ws.Range("K134").Value = 4319.7999
$ws.Range("M134").Value = -1784.7999
$ws.Range("H141").Value = 452780.4
$ws.Range("J141").Value = 452780.4
$ws.Range("L141").Value = 452780.4
$ws.Range("N141").Value = -463140.4

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 187790080
$ws.Range("I4").Value = 139414400
$ws.Range("K4").Value = 418243200
$ws.Range("M4").Value = -418243088
$ws.Range("H68").Value = 573
$ws.Range("I68").Value = 995
$ws.Range("J68").Value = 467.5
$ws.Range("K68").Value = 2985
$ws.Range("L68").Value = 1402.5
$ws.Range("M68").Value = -2174
$ws.Range("N68").Value = -3024.5
$ws.Range("H71").Value = 573
$ws.Range("I71").Value = 995
$ws.Range("J71").Value = 467.5
$ws.Range("K71").Value = 8955
$ws.Range("L71").Value = 4207.5
$ws.Range("M71").Value = -4899
$ws.Range("N71").Value = -12319.5
$ws.Range("H131").Value = 1809.238
$ws.Range("J131").Value = 1991.5333
$ws.Range("L131").Value = 5974.5999
$ws.Range("N131").Value = -16054.5999
$ws.Range("H140").Value = 3054.9167
$ws.Range("I140").Value = 2082.375
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 6247.125
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -1067.125
$ws.Range("N140").Value = -25360

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2874.5833
$ws.Range("I80").Value = 2781.3635
$ws.Range("J80").Value = 3900
$ws.Range("K80").Value = 2781.3635
$ws.Range("L80").Value = 3900
$ws.Range("M80").Value = -1783.3635
$ws.Range("N80").Value = -5896
$ws.Range("H83").Value = 2874.5833
$ws.Range("I83").Value = 2781.3635
$ws.Range("J83").Value = 3900
$ws.Range("K83").Value = 13906.8175
$ws.Range("L83").Value = 19500
$ws.Range("M83").Value = -8914.817499999999
$ws.Range("N83").Value = -29484
$ws.Range("H113").Value = 99498
$ws.Range("I113").Value = 99000
$ws.Range("J113").Value = 99664
$ws.Range("K113").Value = 99000
$ws.Range("L113").Value = 99664
$ws.Range("M113").Value = -96830
$ws.Range("N113").Value = -104004
$ws.Range("H122").Value = 1669
$ws.Range("I122").Value = 1007
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3021
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -571
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 2789.3
$ws.Range("I132").Value = 2570.6667
$ws.Range("K132").Value = 7712.000100000001
$ws.Range("M132").Value = -5182.000100000001

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4363.7896
$ws.Range("I61").Value = 1835.2142
$ws.Range("K61").Value = 1835.2142
$ws.Range("M61").Value = -1633.2142
$ws.Range("H93").Value = 3172.1428
$ws.Range("I93").Value = 2941
$ws.Range("J93").Value = 3750
$ws.Range("K93").Value = 2941
$ws.Range("L93").Value = 3750
$ws.Range("M93").Value = -1693
$ws.Range("N93").Value = -6246
$ws.Range("H113").Value = 4363.7896
$ws.Range("I113").Value = 1835.2142
$ws.Range("K113").Value = 1835.2142
$ws.Range("M113").Value = 334.7858000000001
$ws.Range("H132").Value = 4623.0713
$ws.Range("I132").Value = 4252.3
$ws.Range("K132").Value = 12756.9
$ws.Range("M132").Value = -10226.9
$ws.Range("H136").Value = 7934.8823
$ws.Range("I136").Value = 5659.533
$ws.Range("K136").Value = 16978.599
$ws.Range("M136").Value = -14428.599

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11482.143
$ws.Range("I122").Value = 4875.2
$ws.Range("J122").Value = 27999.5
$ws.Range("K122").Value = 14625.6
$ws.Range("L122").Value = 83998.5
$ws.Range("M122").Value = -12175.6
$ws.Range("N122").Value = -88898.5
$ws.Range("H132").Value = 2334
$ws.Range("I132").Value = 1875.75
$ws.Range("K132").Value = 5627.25
$ws.Range("M132").Value = -3097.25
$ws.Range("H136").Value = 13402193
$ws.Range("J136").Value = 22335956
$ws.Range("L136").Value = 67007868
$ws.Range("N136").Value = -67012968
